$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("customers")

# Howard's birthday changes from 12/05/1987 to 21/05/2002
$ws1.Range("E5").Value = "21/05/2002"

# New customer row: Billy, 23, canDrinkAlcohol=TRUE, Beer, birthday=DATE(1940,4,28), height=1.25
$ws1.Range("A6").Value = "Billy"
$ws1.Range("B6").Value = 23
$ws1.Range("C6").Value = $true
$ws1.Range("C6").NumberFormat = """TRUE"";""TRUE"";""FALSE"""
$ws1.Range("D6").Value = "Beer"
$ws1.Range("E6").NumberFormat = "mm/dd/yy"
$ws1.Range("E6").Formula = "=DATE(1940,4,28)"
$ws1.Range("F6").Value = 1.25

# Activate the customers sheet and set its selection
$ws1.Activate()
$ws1.Range("E7").Select() | Out-Null
